$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data grid (header + 19 data rows) for the FIPE lookup table.
# Column order: MarcaSelecionada, ModeloSelecionado, AnoSelecionado, CodigoFipe, PrecoMedio
$data = @(
  @("Acura","Integra GS 1.8","1992 Gasolina","038003-2"," 11097.00"),
  @("Acura","Integra GS 1.8","1991 Gasolina","038003-2"," 10366.00"),
  @("Acura","Legend 3.2/3.5","1998 Gasolina","038001-6"," 40991.00"),
  @("Acura","Legend 3.2/3.5","1997 Gasolina","038002-4"," 22580.00"),
  @("Acura","Legend 3.2/3.5","1996 Gasolina","038002-4"," 21233.00"),
  @("Acura","Legend 3.2/3.5","1995 Gasolina","038002-4"," 19084.00"),
  @("Acura","Legend 3.2/3.5","1994 Gasolina","038002-4"," 18267.00"),
  @("Acura","Legend 3.2/3.5","1993 Gasolina","038002-4"," 16282.00"),
  @("Acura","Legend 3.2/3.5","1992 Gasolina","038002-4"," 14802.00"),
  @("Acura","Legend 3.2/3.5","1991 Gasolina","038002-4"," 14219.00"),
  @("Acura","NSX 3.0","1995 Gasolina","038001-6"," 40991.00"),
  @("Acura","NSX 3.0","1994 Gasolina","038001-6"," 39550.00"),
  @("Acura","NSX 3.0","1993 Gasolina","038001-6"," 38236.00"),
  @("Acura","NSX 3.0","1992 Gasolina","038001-6"," 36538.00"),
  @("Acura","NSX 3.0","1991 Gasolina","038001-6"," 33397.00"),
  @("Agrale","MARRUÁ 2.8 12V 132cv TDI Diesel","2007 Diesel","006009-7"," 27313.00"),
  @("Agrale","MARRUÁ 2.8 12V 132cv TDI Diesel","2006 Diesel","060001-6"," 44511.00"),
  @("Agrale","MARRUÁ 2.8 12V 132cv TDI Diesel","2005 Diesel","060001-6"," 43362.00"),
  @("Agrale","MARRUÁ 2.8 12V 132cv TDI Diesel","2004 Diesel","060001-6"," 36756.00")
)

# PrecoMedio (column E) values look like numbers ("" 11097.00"") so Excel would
# normally coerce them on assignment. Force that column to text first so the
# values land as plain shared strings, matching the original table's layout.
$priceRange = $ws.Range("E2:E20")
$priceRange.NumberFormat = "@"

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# Drop the temporary text number-format again so the cells fall back to the
# workbook's default (unstyled) look, same as the rest of the data rows.
$priceRange.Style = "Normal"
